# New crime data collected — weekly CompStat report roll-forward
# (Volume/Number string, reporting week dates, and the precinct data table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: "Volume 31   Number  43" -> "...44"; report week 10/21-10/27 -> 10/28-11/3
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# ---------------------------------------------------------------------------
# Column widths for I (9) and J (10) shrink to match the other numeric cols
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 6.168446
$ws.Columns.Item(10).ColumnWidth = 6.168446

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("M14").Value = 100

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("I15").Value = 36
$ws.Range("K15").Value = 63.636363636363
$ws.Range("L15").Value = 2.857142857142
$ws.Range("M15").Value = 28.571428571428
$ws.Range("N15").Value = 33.333333333333

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 55.555555555555
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 371
$ws.Range("J16").Value = 298
$ws.Range("K16").Value = 24.496644295302
$ws.Range("L16").Value = 43.798449612403
$ws.Range("M16").Value = 24.915824915824
$ws.Range("N16").Value = -66.455696202531

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 5.882352941176
$ws.Range("I17").Value = 546
$ws.Range("J17").Value = 407
$ws.Range("K17").Value = 34.152334152334
$ws.Range("L17").Value = 52.089136490250
$ws.Range("M17").Value = 97.826086956521
$ws.Range("N17").Value = 48.773841961852

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "***.*"
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 184
$ws.Range("K18").Value = 50.819672131147
$ws.Range("L18").Value = 46.031746031746
$ws.Range("M18").Value = -27.272727272727
$ws.Range("N18").Value = -89.014925373134

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -4.545454545454
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 82
$ws.Range("H19").Value = -9.756097560975
$ws.Range("I19").Value = 817
$ws.Range("J19").Value = 727
$ws.Range("K19").Value = 12.379642365887
$ws.Range("L19").Value = -1.566265060240
$ws.Range("M19").Value = 94.988066825775
$ws.Range("N19").Value = -33.142389525368

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = 29.629629629629
$ws.Range("I20").Value = 265
$ws.Range("J20").Value = 288
$ws.Range("K20").Value = -7.986111111111
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 39.473684210526
$ws.Range("N20").Value = -86.205101509630

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 18.604651162790
$ws.Range("F21").Value = 196
$ws.Range("G21").Value = 189
$ws.Range("H21").Value = 3.703703703703
$ws.Range("I21").Value = 2223
$ws.Range("J21").Value = 1866
$ws.Range("K21").Value = 19.131832797427
$ws.Range("L21").Value = 18.496801705756
$ws.Range("M21").Value = 51.740614334471
$ws.Range("N21").Value = -64.898152534343

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = 1
$ws.Range("J22").Value = 76
$ws.Range("K22").Value = -14.473684210526
$ws.Range("L22").Value = 16.071428571428

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = -41.860465116279
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 154
$ws.Range("H24").Value = -32.467532467532
$ws.Range("I24").Value = 1718
$ws.Range("J24").Value = 1681
$ws.Range("K24").Value = 2.201070791195
$ws.Range("L24").Value = 6.575682382134
$ws.Range("M24").Value = 79.519331243469

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 40
$ws.Range("H25").Value = -51.807228915662
$ws.Range("I25").Value = 892
$ws.Range("J25").Value = 848
$ws.Range("K25").Value = 5.188679245283
$ws.Range("L25").Value = -1.870187018701

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 28
$ws.Range("E26").Value = -53.571428571428
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = 6.410256410256
$ws.Range("I26").Value = 1128
$ws.Range("J26").Value = 821
$ws.Range("K26").Value = 37.393422655298
$ws.Range("L26").Value = 47.066492829204
$ws.Range("M26").Value = 41.353383458646

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 57
$ws.Range("K27").Value = 39.024390243902
$ws.Range("L27").Value = 23.913043478260

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 150
$ws.Range("F28").Value = 14
$ws.Range("H28").Value = 55.555555555555
$ws.Range("I28").Value = 127
$ws.Range("J28").Value = 129
$ws.Range("K28").Value = -1.550387596899
$ws.Range("L28").Value = 36.559139784946

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("N29").Value = -87.037037037037

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("N30").Value = -90

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------------
$ws.Range("F31").Value = "'0"
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
